$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.927.34'
$ws.Range('E2').Value = '  -1.21%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.810.64'
$ws.Range('E3').Value = '  -0.68%  '
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '310.62'
$ws.Range('E5').Value = '  -0.81%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  +0.20%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4614'
$ws.Range('E7').Value = '  +3.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3710'
$ws.Range('E8').Value = '  -1.42%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07368'
$ws.Range('E9').Value = '  -0.47%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8735'
$ws.Range('E10').Value = '  -0.68%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.47'
$ws.Range('E11').Value = '  -1.88%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.813.83'
$ws.Range('E12').Value = '  -0.42%  '
$ws.Range('E13').Value = '  -1.21%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '92.13'
$ws.Range('E14').Value = '  -0.89%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.523'
$ws.Range('E15').Value = '  -2.86%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.07062'
$ws.Range('E16').Value = '  +0.14%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.002'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008695'
$ws.Range('E18').Value = '  -1.32%  '
$ws.Range('E19').Value = '  +0.10%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.73'
$ws.Range('E20').Value = '  -2.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '26.926.95'
$ws.Range('E21').Value = '  -1.24%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.320'
$ws.Range('E22').Value = '  -0.37%  '
$ws.Range('E23').Value = '  -2.88%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.023.24'
$ws.Range('E24').Value = '  -1.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.896'
$ws.Range('E25').Value = '  -3.35%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '151.43'
$ws.Range('E26').Value = '  +0.29%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.34'
$ws.Range('E27').Value = '  -1.67%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.143'
$ws.Range('E28').Value = '  -6.15%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.303'
$ws.Range('E29').Value = '  -0.83%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '115.88'
$ws.Range('E30').Value = '  -1.29%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08911'
$ws.Range('E31').Value = '  +0.21%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7525'
$ws.Range('E32').Value = '  -4.77%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.155'
$ws.Range('E33').Value = '  -3.58%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.923'
$ws.Range('E34').Value = '  +0.06%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.441'
$ws.Range('E35').Value = '  -2.97%  '
$ws.Range('E36').Value = '  +0.18%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.104'
$ws.Range('E37').Value = '  -0.43%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01978'
$ws.Range('E38').Value = '  +0.22%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05238'
$ws.Range('E39').Value = '  -0.58%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.429'
$ws.Range('E40').Value = '  +3.49%  '
$ws.Range('E41').Value = '  +1.93%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5303'
$ws.Range('E42').Value = '  +0.08%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '7.179'
$ws.Range('E43').Value = '  -1.56%  '
$ws.Range('E44').Value = '  -2.19%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.505'
$ws.Range('E45').Value = '  -1.57%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4978'
$ws.Range('E46').Value = '  -1.57%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.31'
$ws.Range('E47').Value = '  -2.28%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '104.08'
$ws.Range('E48').Value = '  -1.31%  '
$ws.Range('E49').Value = '  +0.17%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.668'
$ws.Range('E50').Value = '  -1.05%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06295'
$ws.Range('E51').Value = '  -1.36%  '
